$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'21.754.23"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.28%  "
$ws.Range("D3").Value = "'1.539.26"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.01%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("E5").Value = "  -0.08%  "
$ws.Range("D6").Value = "'289.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.17%  "
$ws.Range("D7").Value = "'0.3896"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.00%  "
$ws.Range("E8").Value = "  -1.23%  "
$ws.Range("D9").Value = "'43.04"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.31%  "
$ws.Range("D10").Value = "'0.07209"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.26%  "
$ws.Range("D11").Value = "'1.059"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.09%  "
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("D13").Value = "'5.644"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.82%  "
$ws.Range("D14").Value = "'18.62"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.52%  "
$ws.Range("D15").Value = "'6.611"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.52%  "
$ws.Range("D16").Value = "'0.00001115"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.00%  "
$ws.Range("D17").Value = "'1.540.66"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.64%  "
$ws.Range("D18").Value = "'0.06594"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.77%  "
$ws.Range("D19").Value = "'83.20"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.18%  "
$ws.Range("D21").Value = "'6.146"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.67%  "
$ws.Range("D22").Value = "'15.41"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.50%  "
$ws.Range("D23").Value = "'10.91"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.66%  "
$ws.Range("D24").Value = "'2.415"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.10%  "
$ws.Range("D25").Value = "'21.760.20"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.30%  "
$ws.Range("D26").Value = "'2.377"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.68%  "
$ws.Range("D27").Value = "'146.61"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.83%  "
$ws.Range("E28").Value = "  -2.59%  "
$ws.Range("D29").Value = "'4.831"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.29%  "
$ws.Range("D30").Value = "'1.718.18"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.46%  "
$ws.Range("D31").Value = "'117.60"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.13%  "
$ws.Range("D32").Value = "'0.9703"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -12.14%  "
$ws.Range("D33").Value = "'5.919"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.45%  "
$ws.Range("D34").Value = "'0.08185"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.19%  "
$ws.Range("D35").Value = "'8.900"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.86%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "'0.06077"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.74%  "
$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").Value = "'5.139"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.83%  "
$ws.Range("D38").Value = "'1.493"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.07%  "
$ws.Range("D39").Value = "'0.02210"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.01%  "
$ws.Range("D40").Value = "'0.2039"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.84%  "
$ws.Range("D41").Value = "'1.192"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.42%  "
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("D43").Value = "'10.68"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.00%  "
$ws.Range("D44").Value = "'0.5755"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.05%  "
$ws.Range("D45").Value = "'13.07"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.16%  "
$ws.Range("D46").Value = "'3.745"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.75%  "
$ws.Range("D47").Value = "'0.5516"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.84%  "
$ws.Range("D48").Value = "'117.55"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.66%  "
$ws.Range("B49").Value = "EOS"
$ws.Range("C49").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D49").Value = "'1.161"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.53%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "'1.868"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.76%  "
$ws.Range("D51").Value = "'0.06729"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.92%  "
